# Apply genitive-case wording/name updates to the "print_docs" signature
# block of the template, per commit "genitive case in print_docs func is
# added".

$d = $word.ActiveDocument

# 1) Job title + signer name -> new signer (genitive-ish phrasing change)
$d.Content.Find.Execute(
    "специалист 1-ой категории Петряков О. Ю.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "генерального директора Котлярчука О. Ю.", 2)

# 2) "Доверенность" -> genitive "Доверенности"
$d.Content.Find.Execute(
    "Доверенность № 123456 от 02.02.2024", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Доверенности № 123456 от 02.02.2024", 2)

# 3) Fix double space before signer's surname
$d.Content.Find.Execute(
    "В. Г.  Кемоклидзе", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "В. Г. Кемоклидзе", 2)

# 4) Printed signature block: surname change Петряков -> Котлярчук
$d.Content.Find.Execute(
    "О. Ю. Петряков", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "О. Ю. Котлярчук", 2)
